$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: date changed
$ws.Range("D2").Value = 44559

# Row 3: date changed
$ws.Range("D3").Value = 44559

# Row 4: date and volume/price values changed
$ws.Range("D4").Value = 44574
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6500
$ws.Range("S4").Value = 3250

# Row 5: date and volume/price values changed
$ws.Range("D5").Value = 44574
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("S5").Value = 2500

# Row 6: date and volume/price values changed
$ws.Range("D6").Value = 44223
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 3500
$ws.Range("O6").Value = 4000
$ws.Range("P6").Value = 3750
$ws.Range("S6").Value = 1875

# Row 7: date and volume/price values changed
$ws.Range("D7").Value = 44223
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 3000
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 3000
$ws.Range("S7").Value = 1500
